# Generate Report for Handoff
#
# A new handoff run completed for the "xinjiang" localization job: the
# source markdown file was re-uploaded under a new GUID and new handoff
# (.xlf) packages were produced (new GUID + new content hash) for the
# zh-cn and de-de targets. Refresh the status report with the new file
# names / hashes and the new handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "fb950e04-a55d-4146-8067-eabfb4028a28"
$newGuid = "f791a250-c264-4e1d-bbb2-e2b94ff1eadb"

$oldHash = "0d9f5ddc83662c13605f56d4926f44caea771944"
$newHash = "0239ecc090eae3b8f503f704533a1bc401a9d926"

$oldMdName = "$oldGuid.md"
$newMdName = "$newGuid.md"

$oldZhXlf = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"

$oldDeXlf = "$oldGuid.$oldHash.de-de.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$oldOverviewDate = "2016-03-25 08:39:52"
$newOverviewDate = "2016-03-25 08:40:53"

$oldZhDate = "2016-03-25 08:39:43"
$newZhDate = "2016-03-25 08:40:44"

function Update-DisplayText($ws, $oldText, $newText) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.TextToDisplay -eq $oldText) {
            $h.TextToDisplay = $newText
        }
    }
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newOverviewDate
Update-DisplayText $wsOverview $oldMdName $newMdName

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("D2").Value = $newZhXlf
$wsZhCn.Range("E2").Value = $newZhDate
Update-DisplayText $wsZhCn $oldMdName $newMdName
Update-DisplayText $wsZhCn $oldZhXlf $newZhXlf

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("D2").Value = $newDeXlf
$wsDeDe.Range("E2").Value = $newOverviewDate
Update-DisplayText $wsDeDe $oldMdName $newMdName
Update-DisplayText $wsDeDe $oldDeXlf $newDeXlf
